$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.160.16"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.575.96"
$ws.Range("E3").Value = "  -2.76%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "588.14"
$ws.Range("E5").Value = "  -3.37%  "

# Row 6 - Solana
$ws.Range("D6").Value = "149.45"
$ws.Range("E6").Value = "  +0.83%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.65%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.110"
$ws.Range("E9").Value = "  +0.34%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +1.46%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.383"
$ws.Range("E11").Value = "  +0.22%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.50%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "27.35"
$ws.Range("E13").Value = "  -0.56%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.039.45"
$ws.Range("E14").Value = "  -2.57%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "63.047.23"
$ws.Range("E15").Value = "  -0.81%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +5.36%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.596.04"
$ws.Range("E17").Value = "  -2.57%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "12.30"
$ws.Range("E18").Value = "  +4.79%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "4.72"
$ws.Range("E19").Value = "  +3.10%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "345.11"
$ws.Range("E20").Value = "  -0.43%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.87"
$ws.Range("E21").Value = "  -0.46%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.15%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "67.37"
$ws.Range("E23").Value = "  +1.43%  "

# Row 24 - SuiNetwork
$ws.Range("E24").Value = "  +2.83%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = "9.22"
$ws.Range("E25").Value = "  +0.91%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  -1.36%  "

# Row 27 - Bittensor
$ws.Range("D27").Value = "550.99"
$ws.Range("E27").Value = "  -3.00%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.70%  "

# Row 29 - Aptos
$ws.Range("E29").Value = "  -1.65%  "

# Row 30 - Kaspa
$ws.Range("D30").Value = "0.161"
$ws.Range("E30").Value = "  -1.06%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.71%  "

# Row 32 - PEPE
$ws.Range("D32").Value = "0.0$([char]0x2083)0844"
$ws.Range("E32").Value = "  -1.20%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  -1.28%  "

# Row 34 - NEARProtocol
$ws.Range("D34").Value = "5.18"
$ws.Range("E34").Value = "  -2.27%  "

# Row 35 - Monero
$ws.Range("D35").Value = "167.50"
$ws.Range("E35").Value = "  -0.87%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "  +1.31%  "

# Row 37 - FirstDigitalUSD
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.03%  "

# Row 38 - EthereumClassic
$ws.Range("D38").Value = "19.47"
$ws.Range("E38").Value = "  +1.59%  "

# Row 39 - Stacks
$ws.Range("D39").Value = "1.92"
$ws.Range("E39").Value = "  -0.65%  "

# Row 40 - USDe
$ws.Range("E40").Value = "  +0.08%  "

# Row 41 - Aave
$ws.Range("D41").Value = "166.02"
$ws.Range("E41").Value = "  +0.36%  "

# Row 42 - OKB
$ws.Range("D42").Value = "39.54"
$ws.Range("E42").Value = "  -1.25%  "

# Row 43 - Filecoin
$ws.Range("D43").Value = "3.93"
$ws.Range("E43").Value = "  +3.48%  "

# Row 44 & 45 - InjectiveProtocol and Hedera swap places
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "0.0581"
$ws.Range("E44").Value = "  +2.16%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "22.34"
$ws.Range("E45").Value = "  +1.26%  "

# Row 46 - Mantle
$ws.Range("D46").Value = "0.627"
$ws.Range("E46").Value = "  -0.54%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  +2.34%  "

# Row 48 - dogwifhat
$ws.Range("E48").Value = "  +0.36%  "

# Row 49 - Stellar
$ws.Range("D49").Value = "0.0961"

# Row 50 - EnergySwap
$ws.Range("D50").Value = "18.96"
$ws.Range("E50").Value = "  +0.19%  "

# Row 51 - BabyDogeCoin
$ws.Range("D51").Value = "0.0$([char]0x2086)0230"
$ws.Range("E51").Value = "  +16.81%  "
